$d = $word.ActiveDocument

# Locate the paragraph "Opção de participar num livechat com o vendedor."
# (last bullet of the "Nice to Have" list) so we can append a new bullet
# right after it: "Deployment da aplicação."
$target = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "*participar num livechat com o vendedor*") {
        $target = $cand
    }
}

if ($target -eq $null) {
    Write-Output "ERROR: target paragraph not found"
} else {
    $rng = $target.Range
    # Collapse to the end of the paragraph's text (just before its
    # paragraph mark) and insert a literal carriage return followed by
    # the new sentence. Typing a CR this way splits the paragraph the
    # same way pressing Enter would, so the brand-new paragraph/run
    # inherit the exact pPr/rPr (list style, numbering, spacing, fonts,
    # color, size, language) of the paragraph it was split from.
    $rng.Collapse(0)
    $cr = [char]13
    $rng.InsertAfter("$($cr)Deployment da aplicação.")
    Write-Output "Inserted new bullet after target paragraph."
}
